$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Melbourne Victory vs Wellington Phoenix)
$ws.Range("F2").Value = 1.69
$ws.Range("H2").Value = 5.3
$ws.Range("I2").Value = 5.5
$ws.Range("V2").Value = 1.22
$ws.Range("W2").Value = 2.4
$ws.Range("Z2").Value = 50

# Row 5 (MC Oran vs Belouizdad)
$ws.Range("J5").Value = 2.72
$ws.Range("P5").Value = 1.42
$ws.Range("Q5").Value = 2.96

# Row 6 (Al-Taawoun Buraidah vs Al Najma Club)
$ws.Range("F6").Value = 1.47
$ws.Range("G6").Value = 1.73
$ws.Range("I6").Value = 17.5
$ws.Range("K6").Value = 8.800000000000001

# Row 10 (Porto vs AVS Futebol SAD)
$ws.Range("H10").Value = 30
$ws.Range("I10").Value = 130
$ws.Range("J10").Value = 11
